# Update res_bus/vm_pu.xlsx voltage-magnitude results for the
# "case with 380 kV" re-run (slack bus vm_pu dropped 1.05 -> 1.02 p.u.,
# which cascades into new load-flow solutions across buses B..N).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.033402712836467
$ws.Cells.Item(2, 4).Value = 1.035755606207742
$ws.Cells.Item(2, 5).Value = 1.052201421705744
$ws.Cells.Item(2, 6).Value = 1.057679932864867
$ws.Cells.Item(2, 9).Value = 1.035777505714085
$ws.Cells.Item(2, 10).Value = 1.038527633493211
$ws.Cells.Item(2, 11).Value = 1.03855130136421
$ws.Cells.Item(2, 12).Value = 1.054950809217675
$ws.Cells.Item(2, 13).Value = 1.060414240607005
$ws.Cells.Item(2, 14).Value = 1.016798842042896
# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.034329766068087
$ws.Cells.Item(3, 4).Value = 1.036427912282907
$ws.Cells.Item(3, 5).Value = 1.053467042264478
$ws.Cells.Item(3, 6).Value = 1.059016908378663
$ws.Cells.Item(3, 9).Value = 1.03598388010999
$ws.Cells.Item(3, 10).Value = 1.039097499322869
$ws.Cells.Item(3, 11).Value = 1.039033432436287
$ws.Cells.Item(3, 12).Value = 1.056028065862906
$ws.Cells.Item(3, 13).Value = 1.061563772226247
$ws.Cells.Item(3, 14).Value = 1.016989481453113
# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.034929596865193
$ws.Cells.Item(4, 4).Value = 1.036862901629136
$ws.Cells.Item(4, 5).Value = 1.054286710371128
$ws.Cells.Item(4, 6).Value = 1.059882788894382
$ws.Cells.Item(4, 9).Value = 1.036116226252016
$ws.Cells.Item(4, 10).Value = 1.03946558768959
$ws.Cells.Item(4, 11).Value = 1.039344686094113
$ws.Cells.Item(4, 12).Value = 1.056725292796357
$ws.Cells.Item(4, 13).Value = 1.062307819885126
$ws.Cells.Item(4, 14).Value = 1.017112571294247
# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.035181757080098
$ws.Cells.Item(5, 4).Value = 1.037045761248684
$ws.Cells.Item(5, 5).Value = 1.054631474390235
$ws.Cells.Item(5, 6).Value = 1.060246990468926
$ws.Cells.Item(5, 9).Value = 1.036171578949481
$ws.Cells.Item(5, 10).Value = 1.039620175332636
$ws.Cells.Item(5, 11).Value = 1.039475364888434
$ws.Cells.Item(5, 12).Value = 1.057018448953602
$ws.Cells.Item(5, 13).Value = 1.06262067237744
$ws.Cells.Item(5, 14).Value = 1.01716425424827
# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.035224095371918
$ws.Cells.Item(6, 4).Value = 1.037076463576849
$ws.Cells.Item(6, 5).Value = 1.054689372124217
$ws.Cells.Item(6, 6).Value = 1.060308152433503
$ws.Cells.Item(6, 9).Value = 1.036180856153912
$ws.Cells.Item(6, 10).Value = 1.039646122077351
$ws.Cells.Item(6, 11).Value = 1.039497296318025
$ws.Cells.Item(6, 12).Value = 1.057067673660271
$ws.Cells.Item(6, 13).Value = 1.062673204941389
$ws.Cells.Item(6, 14).Value = 1.017172928284268
# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.03493296627677
$ws.Cells.Item(7, 4).Value = 1.036865345048353
$ws.Cells.Item(7, 5).Value = 1.054291316433145
$ws.Cells.Item(7, 6).Value = 1.059887654643168
$ws.Cells.Item(7, 9).Value = 1.036116967000109
$ws.Cells.Item(7, 10).Value = 1.039467653913485
$ws.Cells.Item(7, 11).Value = 1.039346432907987
$ws.Cells.Item(7, 12).Value = 1.0567292097968
$ws.Cells.Item(7, 13).Value = 1.062312000018293
$ws.Cells.Item(7, 14).Value = 1.017113262136183
# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.033716021770947
$ws.Cells.Item(8, 4).Value = 1.035982822452979
$ws.Cells.Item(8, 5).Value = 1.052628994479679
$ws.Cells.Item(8, 6).Value = 1.058131611922911
$ws.Cells.Item(8, 9).Value = 1.035847497528594
$ws.Cells.Item(8, 10).Value = 1.038720357276022
$ws.Cells.Item(8, 11).Value = 1.038714388056568
$ws.Cells.Item(8, 12).Value = 1.055314838752053
$ws.Cells.Item(8, 13).Value = 1.060802684616251
$ws.Cells.Item(8, 14).Value = 1.016863324612869
# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.031571360328101
$ws.Cells.Item(9, 4).Value = 1.034427455600163
$ws.Cells.Item(9, 5).Value = 1.049705271140779
$ws.Cells.Item(9, 6).Value = 1.055043052356155
$ws.Cells.Item(9, 9).Value = 1.035363541631395
$ws.Cells.Item(9, 10).Value = 1.037398541499922
$ws.Cells.Item(9, 11).Value = 1.037595171767733
$ws.Cells.Item(9, 12).Value = 1.052823779307446
$ws.Cells.Item(9, 13).Value = 1.05814473548057
$ws.Cells.Item(9, 14).Value = 1.016420869130226
# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.030141436531699
$ws.Cells.Item(10, 4).Value = 1.033390427965047
$ws.Cells.Item(10, 5).Value = 1.047759724202017
$ws.Cells.Item(10, 6).Value = 1.052987821071018
$ws.Cells.Item(10, 9).Value = 1.03503479077611
$ws.Cells.Item(10, 10).Value = 1.036514001570785
$ws.Cells.Item(10, 11).Value = 1.036845376216197
$ws.Cells.Item(10, 12).Value = 1.051163826205372
$ws.Cells.Item(10, 13).Value = 1.056373802297804
$ws.Cells.Item(10, 14).Value = 1.016124542064236
# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.029522231342805
$ws.Cells.Item(11, 4).Value = 1.032941366299381
$ws.Cells.Item(11, 5).Value = 1.046918113234335
$ws.Cells.Item(11, 6).Value = 1.052098763803399
$ws.Cells.Item(11, 9).Value = 1.034890990994194
$ws.Cells.Item(11, 10).Value = 1.03613020025647
$ws.Cells.Item(11, 11).Value = 1.0365198458784
$ws.Cells.Item(11, 12).Value = 1.050445207899062
$ws.Cells.Item(11, 13).Value = 1.055607194910095
$ws.Cells.Item(11, 14).Value = 1.015995909414484
# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.029292225289991
$ws.Cells.Item(12, 4).Value = 1.032774562301759
$ws.Cells.Item(12, 5).Value = 1.046605623225413
$ws.Cells.Item(12, 6).Value = 1.051768657230261
$ws.Cells.Item(12, 9).Value = 1.034837359930429
$ws.Cells.Item(12, 10).Value = 1.035987520960754
$ws.Cells.Item(12, 11).Value = 1.036398800002729
$ws.Cells.Item(12, 12).Value = 1.050178302241944
$ws.Cells.Item(12, 13).Value = 1.055322473898689
$ws.Cells.Item(12, 14).Value = 1.015948081423387
# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.029341562623578
$ws.Cells.Item(13, 4).Value = 1.032810342426101
$ws.Cells.Item(13, 5).Value = 1.04667264791409
$ws.Cells.Item(13, 6).Value = 1.051839460411394
$ws.Cells.Item(13, 9).Value = 1.034848873811293
$ws.Cells.Item(13, 10).Value = 1.03601813150793
$ws.Cells.Item(13, 11).Value = 1.036424770604138
$ws.Cells.Item(13, 12).Value = 1.050235553480479
$ws.Cells.Item(13, 13).Value = 1.055383546141843
$ws.Cells.Item(13, 14).Value = 1.01595834286411
# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.029503219085662
$ws.Cells.Item(14, 4).Value = 1.032927578275398
$ws.Cells.Item(14, 5).Value = 1.046892280239902
$ws.Cells.Item(14, 6).Value = 1.052071474473749
$ws.Cells.Item(14, 9).Value = 1.034886562267933
$ws.Cells.Item(14, 10).Value = 1.036118408755276
$ws.Cells.Item(14, 11).Value = 1.036509842826372
$ws.Cells.Item(14, 12).Value = 1.050423144972025
$ws.Cells.Item(14, 13).Value = 1.055583659157759
$ws.Cells.Item(14, 14).Value = 1.015991956918189
# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.029602820219602
$ws.Cells.Item(15, 4).Value = 1.032999810843436
$ws.Cells.Item(15, 5).Value = 1.047027619042295
$ws.Cells.Item(15, 6).Value = 1.052214443011977
$ws.Cells.Item(15, 9).Value = 1.034909754565328
$ws.Cells.Item(15, 10).Value = 1.036180177195189
$ws.Cells.Item(15, 11).Value = 1.036562241492953
$ws.Cells.Item(15, 12).Value = 1.050538729036534
$ws.Cells.Item(15, 13).Value = 1.05570695946238
$ws.Cells.Item(15, 14).Value = 1.016012661277091
# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.030182530336846
$ws.Cells.Item(16, 4).Value = 1.033420230315301
$ws.Cells.Item(16, 5).Value = 1.04781559633594
$ws.Cells.Item(16, 6).Value = 1.053046843059651
$ws.Cells.Item(16, 9).Value = 1.035044303792824
$ws.Cells.Item(16, 10).Value = 1.036539456582933
$ws.Cells.Item(16, 11).Value = 1.036866962423332
$ws.Cells.Item(16, 12).Value = 1.051211521599335
$ws.Cells.Item(16, 13).Value = 1.05642468396079
$ws.Cells.Item(16, 14).Value = 1.016133072236875
# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.030546156906655
$ws.Cells.Item(17, 4).Value = 1.033683943045326
$ws.Cells.Item(17, 5).Value = 1.04831009280709
$ws.Cells.Item(17, 6).Value = 1.053569217484833
$ws.Cells.Item(17, 9).Value = 1.035128315305342
$ws.Cells.Item(17, 10).Value = 1.03676461180085
$ws.Cells.Item(17, 11).Value = 1.037057874786087
$ws.Cells.Item(17, 12).Value = 1.051633586177934
$ws.Cells.Item(17, 13).Value = 1.056874950837896
$ws.Cells.Item(17, 14).Value = 1.016208516989236
# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.030758250384255
$ws.Cells.Item(18, 4).Value = 1.033837760057981
$ws.Cells.Item(18, 5).Value = 1.048598603839393
$ws.Cells.Item(18, 6).Value = 1.053873993867472
$ws.Cells.Item(18, 9).Value = 1.035177178027475
$ws.Cells.Item(18, 10).Value = 1.036895864820484
$ws.Cells.Item(18, 11).Value = 1.037169147356561
$ws.Cells.Item(18, 12).Value = 1.05187978411902
$ws.Cells.Item(18, 13).Value = 1.057137605072583
$ws.Cells.Item(18, 14).Value = 1.016252491648459
# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.030830568138872
$ws.Cells.Item(19, 4).Value = 1.033890207293857
$ws.Cells.Item(19, 5).Value = 1.048696992208407
$ws.Cells.Item(19, 6).Value = 1.053977929084897
$ws.Cells.Item(19, 9).Value = 1.03519381523928
$ws.Cells.Item(19, 10).Value = 1.03694060577441
$ws.Cells.Item(19, 11).Value = 1.037207074273156
$ws.Cells.Item(19, 12).Value = 1.051963733819212
$ws.Cells.Item(19, 13).Value = 1.057227167023861
$ws.Cells.Item(19, 14).Value = 1.016267480610395
# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.030507143603608
$ws.Cells.Item(20, 4).Value = 1.033655649382963
$ws.Cells.Item(20, 5).Value = 1.0482570297955
$ws.Cells.Item(20, 6).Value = 1.053513162956118
$ws.Cells.Item(20, 9).Value = 1.035119316120914
$ws.Cells.Item(20, 10).Value = 1.036740462657709
$ws.Cells.Item(20, 11).Value = 1.037037400318
$ws.Cells.Item(20, 12).Value = 1.051588301122964
$ws.Cells.Item(20, 13).Value = 1.056826639284254
$ws.Cells.Item(20, 14).Value = 1.016200425685549
# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.029455615442258
$ws.Cells.Item(21, 4).Value = 1.032893055302564
$ws.Cells.Item(21, 5).Value = 1.046827600648331
$ws.Cells.Item(21, 6).Value = 1.052003148573819
$ws.Cells.Item(21, 9).Value = 1.034875469957138
$ws.Cells.Item(21, 10).Value = 1.036088882871719
$ws.Cells.Item(21, 11).Value = 1.036484794746711
$ws.Cells.Item(21, 12).Value = 1.050367903382445
$ws.Cells.Item(21, 13).Value = 1.055524730017239
$ws.Cells.Item(21, 14).Value = 1.015982059743678
# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.028794445187818
$ws.Cells.Item(22, 4).Value = 1.032413568054482
$ws.Cells.Item(22, 5).Value = 1.04592956539998
$ws.Cells.Item(22, 6).Value = 1.051054487257931
$ws.Cells.Item(22, 9).Value = 1.034720896356798
$ws.Cells.Item(22, 10).Value = 1.035678523899578
$ws.Cells.Item(22, 11).Value = 1.036136601128671
$ws.Cells.Item(22, 12).Value = 1.04960071257292
$ws.Cells.Item(22, 13).Value = 1.054706347115675
$ws.Cells.Item(22, 14).Value = 1.015844486199956
# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.029144947106325
$ws.Cells.Item(23, 4).Value = 1.03266775439132
$ws.Cells.Item(23, 5).Value = 1.046405564842769
$ws.Cells.Item(23, 6).Value = 1.051557320667251
$ws.Cells.Item(23, 9).Value = 1.034802957909509
$ws.Cells.Item(23, 10).Value = 1.03589612779473
$ws.Cells.Item(23, 11).Value = 1.036321255982384
$ws.Cells.Item(23, 12).Value = 1.050007403972298
$ws.Cells.Item(23, 13).Value = 1.055140170964448
$ws.Cells.Item(23, 14).Value = 1.015917442862891
# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.030524772047426
$ws.Cells.Item(24, 4).Value = 1.033668434077809
$ws.Cells.Item(24, 5).Value = 1.048281006439873
$ws.Cells.Item(24, 6).Value = 1.05353849132404
$ws.Cells.Item(24, 9).Value = 1.035123382896844
$ws.Cells.Item(24, 10).Value = 1.036751374851425
$ws.Cells.Item(24, 11).Value = 1.037046652106425
$ws.Cells.Item(24, 12).Value = 1.051608763443799
$ws.Cells.Item(24, 13).Value = 1.056848469129873
$ws.Cells.Item(24, 14).Value = 1.016204081893137
# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.032125834415625
$ws.Cells.Item(25, 4).Value = 1.034829579463566
$ws.Cells.Item(25, 5).Value = 1.050460481804124
$ws.Cells.Item(25, 6).Value = 1.05584084051551
$ws.Cells.Item(25, 9).Value = 1.035489734541382
$ws.Cells.Item(25, 10).Value = 1.037740850947302
$ws.Cells.Item(25, 11).Value = 1.037885161629396
$ws.Cells.Item(25, 12).Value = 1.05346763909274
$ws.Cells.Item(25, 13).Value = 1.058831691030374
$ws.Cells.Item(25, 14).Value = 1.016535494382707
